# Workbook under edit
$wb = $excel.ActiveWorkbook

# --- AddCostemerTest sheet -------------------------------------------------
# Update the "runmode" value in row 3 from "N" to "y", and leave the
# selection on E4.
$wsAddCustomer = $wb.Worksheets.Item("AddCostemerTest")
$wsAddCustomer.Activate()
$wsAddCustomer.Range("E3").Value = "y"
$wsAddCustomer.Range("E4").Select()

# --- OpenAccountTest sheet --------------------------------------------------
# Move the saved selection from A3 to F24.
$wsOpenAccount = $wb.Worksheets.Item("OpenAccountTest")
$wsOpenAccount.Activate()
$wsOpenAccount.Range("F24").Select()

# --- test_suite sheet -------------------------------------------------------
# Make this the active (selected) tab when the workbook is saved.
$wsTestSuite = $wb.Worksheets.Item("test_suite")
$wsTestSuite.Activate()
